# Delete the "Бібліотека Keras" slide (slide 2 in the deck).
# All following slides shift up by one position; their content is
# unchanged. This matches the author's commit "Add files via upload"
# which, per the canonical OOXML diff, removes the sldId pointing at
# the Keras-library slide while leaving every other slide's content
# identical.
$p = $ppt.ActivePresentation
$p.Slides.Item(2).Delete()
